$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update row 2
$ws.Range("B2").Value = "Loja Kings Mauá"
$ws.Range("D2").Value = "(11) 4541-5875"

# Update row 3
$ws.Range("B3").Value = "JS Calçados & Roupas"
$ws.Range("C3").Value = "Loja de calçado"
$ws.Range("D3").Value = ""
$ws.Range("D3").Style = "Normal"

# Update row 4
$ws.Range("B4").Value = "Rosi calçados"
$ws.Range("C4").Value = "Loja de calçado"
$ws.Range("D4").Value = "(11) 94310-5100"

# Delete rows 5 and 6 (shift cells up)
$ws.Range("A5:D6").Delete()
